$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 61: 2025/10/05, 日 (Sunday), hour 0, ranking 5
# The date-like text must stay as plain text (matching the other "日付"
# column entries), so enter it with a leading apostrophe to suppress
# Excel's automatic date recognition, then clear the resulting
# quote-prefix formatting so the cell keeps the workbook's default style.
$ws.Range("A61").Value = "'2025/10/05"
$ws.Range("A61").ClearFormats()

$ws.Range("B61").Value = "日"
$ws.Range("C61").Value = 0
$ws.Range("D61").Value = 5
